# Weekly refresh: insert one new daily price record for
# "Feria Lagunitas de Puerto Montt" - Frutilla, pushing the existing
# rows 170-184 down to 171-185 and populating the new row 170 with the
# latest observation (2022-01-17 / serial 44578).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 170:184 down to 171:185 by inserting a new
# row at 170 (mirrors Excel's Rows.Insert, which also carries the
# existing formatting - e.g. the date style on column D - down with it).
$ws.Rows.Item(170).Insert()

# Populate the newly inserted row with this week's observation.
$ws.Range("A170").Value = 4
$ws.Range("B170").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C170").Value = "Los Lagos"
$ws.Range("D170").Value = 44578
$ws.Range("E170").Value = 10
$ws.Range("F170").Value = "Fruta"
$ws.Range("G170").Value = 100101
$ws.Range("H170").Value = "Berries"
$ws.Range("I170").Value = 100112025
$ws.Range("J170").Value = "Frutilla"
$ws.Range("K170").Value = "Sin especificar"
$ws.Range("L170").Value = "Primera"
$ws.Range("M170").Value = 500
$ws.Range("N170").Value = 8500
$ws.Range("O170").Value = 9000
$ws.Range("P170").Value = 8750
$ws.Range("Q170").Value = "`$/caja 7 kilos"
$ws.Range("R170").Value = "Región de La Araucanía"
$ws.Range("S170").Value = 1250
$ws.Range("T170").Value = 7
